$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 341; this pushes the existing rows 341-356
# down to 342-357 (their contents, incl. styles, move with them).
$ws.Rows.Item(341).Insert()

# Populate the newly inserted row 341 with the new record (a fresh
# weekly price observation prepended ahead of the previous ones).
$ws.Cells.Item(341, 1).Value = 10
$ws.Cells.Item(341, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(341, 3).Value = "La Araucanía"
$ws.Cells.Item(341, 4).Value = 45008
$ws.Cells.Item(341, 5).Value = 9
$ws.Cells.Item(341, 6).Value = 100112039
$ws.Cells.Item(341, 7).Value = "Ciboulette"
$ws.Cells.Item(341, 8).Value = "Sin especificar"
$ws.Cells.Item(341, 9).Value = "Primera"
$ws.Cells.Item(341, 10).Value = 65
$ws.Cells.Item(341, 11).Value = 5000
$ws.Cells.Item(341, 12).Value = 5000
$ws.Cells.Item(341, 13).Value = 5000
$ws.Cells.Item(341, 14).Value = "$/docena de atados"
$ws.Cells.Item(341, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(341, 16).Value = 1667
$ws.Cells.Item(341, 17).Value = 3
$ws.Cells.Item(341, 18).Value = "Hortaliza"
